# Regression analysis: resAll fixed
# Adds a "model" column (D) to the features sheet recording which variables
# were kept in the final regression model, moves/adds rationale text in
# column K for excluded variables, fixes up a couple of uni/description
# cells, re-points the AutoFilter at the new "model" column, and updates
# the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "model" (column D) values for rows that didn't have one yet ---
$ws.Range("D2").Value = "no"
$ws.Range("D3").Value = "no"
$ws.Range("D4").Value = "no"
$ws.Range("D5").Value = "no"
$ws.Range("D6").Value = "no"
$ws.Range("D7").Value = "no"
$ws.Range("D8").Value = "no"
$ws.Range("D9").Value = "no"
$ws.Range("D10").Value = "no"
$ws.Range("D11").Value = "no"
$ws.Range("D12").Value = "no"
$ws.Range("D13").Value = "no"
$ws.Range("D14").Value = "no"
$ws.Range("D16").Value = "no"

# --- Flip existing "model" values that changed ---
$ws.Range("D17").Value = "yes"
$ws.Range("D18").Value = "yes"
$ws.Range("D19").Value = "yes"
$ws.Range("D20").Value = "yes"
$ws.Range("D21").Value = "yes"
$ws.Range("D22").Value = "no"

# --- Row 25 (studio): add model flag, rewrite Description, move the old
#     Description text into the new Rationale column ---
$ws.Range("D25").Value = "no"
$ws.Range("J25").Value = "The studio that produced the film"
$ws.Range("K25").Value = "Not a variable that Paramount can change"

$ws.Range("D26").Value = "yes"

# --- Row 29 (top200_box): now flagged redundant with box office success ---
$ws.Range("B29").Value = "no"
$ws.Range("K29").Value = "Redundant with box office success."

$ws.Range("D30").Value = "yes"

# --- Box office rows: not in the data set ---
$ws.Range("D31").Value = "no"
$ws.Range("K31").Value = "Not in data set"
$ws.Range("D32").Value = "no"
$ws.Range("K32").Value = "Not in data set"

# --- Redundant-with-*_log rationale rows ---
$ws.Range("D33").Value = "no"
$ws.Range("K33").Value = "Redundant with cast_experience_log"

$ws.Range("D35").Value = "no"
$ws.Range("K35").Value = "Redundant with cast_votes_log"

$ws.Range("D37").Value = "yes"

$ws.Range("D38").Value = "no"
$ws.Range("K38").Value = "Redundant with director_experience_log"

$ws.Range("K40").Value = "Redundant with imdb_num_votes_log"

$ws.Range("D41").Value = "no"
$ws.Range("K41").Value = "Response variable"

$ws.Range("D42").Value = "yes"

$ws.Range("D43").Value = "no"
$ws.Range("K43").Value = "Redundant"

$ws.Range("D44").Value = "yes"
$ws.Range("D45").Value = "yes"

$ws.Range("D46").Value = "no"
$ws.Range("K46").Value = "Redundant"

$ws.Range("D47").Value = "no"
$ws.Range("K47").Value = "Redundant"

$ws.Range("D48").Value = "yes"

# --- Re-point the AutoFilter at the new "model" column (D, colId 3) showing
#     only "yes", replacing the old uni (colId 1) + Categorical (colId 5)
#     filters. Row visibility is recomputed automatically from this. ---
$ws.AutoFilterMode = $false
$rng = $ws.Range("A1:K48")
[void]$rng.AutoFilter(4, @("yes"))

# --- Restore the saved selection ---
[void]$ws.Range("I48").Select()
